{"js": "// Rebuild the \"Requisitos\" (course prerequisites) bullet list in its new\n// order, with three entries dropped, four entries added, and a typo\n// (\"\u00c0lgebra\" -> \"\u00c1lgebra\") fixed along the way.\nconst newItems = [\n  \"LOB1268 -  Leitura, Escrita e Comunica\u00e7\u00e3o Cient\u00edfica  (Requisito)\",\n  \"LOB1270 -  Qu\u00edmica Experimental Aplicada  (Requisito)\",\n  \"LOM3081 -  Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos S\u00f3lidos  (Requisito)\",\n  \"LOQ4097 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito)\",\n  \"LOQ4098 -  Fundamentos de Qu\u00edmica para Engenharia II (Requisito)\",\n  \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n  \"LOB1040 -  Laborat\u00f3rio de Eletricidade  (Requisito)\",\n  \"LOB1053 -  F\u00edsica III  (Requisito)\",\n  \"LOB1056 -  Introdu\u00e7\u00e3o aos M\u00e9todos Num\u00e9ricos e Computacionais  (Requisito)\",\n  \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n  \"LOB1011 -  Eletricidade Aplicada  (Requisito)\",\n  \"LOB1012 -  Estat\u00edstica  (Requisito)\",\n  \"LOB1024 -  Mec\u00e2nica  (Requisito)\",\n  \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n  \"LOB1037 -  \u00c1lgebra Linear  (Requisito)\",\n  \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n  \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n  \"LOB1041 -  F\u00edsica Experimental II  (Requisito)\",\n  \"LOB1042 -  F\u00edsica Experimental IV  (Requisito)\",\n  \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n  \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n  \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n  \"LOB1009 -  Leitura e Interpreta\u00e7\u00e3o de Desenho T\u00e9cnico  (Requisito)\",\n  \"LOB1018 -  F\u00edsica I  (Requisito)\",\n  \"LOB1019 -  F\u00edsica II  (Requisito)\",\n  \"LOB1021 -  F\u00edsica IV  (Requisito)\",\n  \"LOQ4233 -  Gest\u00e3o de Neg\u00f3cios  (Requisito)\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the \"Requisitos\" heading; the list itself is the very next paragraph.\nlet listParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Requisitos\") {\n    listParagraph = paragraphs.items[i + 1];\n    break;\n  }\n}\nif (!listParagraph) {\n  throw new Error(\"Could not find the Requisitos list paragraph\");\n}\n\n// Clear the existing bullet items (each was \"Code - Name  (Requisito)\"\n// followed by a manual line break).\nlistParagraph.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-insert every item in the new order, each followed by a line break,\n// matching the original \"one line per requirement\" formatting.\nfor (const item of newItems) {\n  const end = listParagraph.getRange(Word.RangeLocation.end);\n  end.insertText(item, Word.InsertLocation.end);\n  await context.sync();\n\n  const afterText = listParagraph.getRange(Word.RangeLocation.end);\n  afterText.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\n  await context.sync();\n}\n", "ps1": "# Rebuild the \"Requisitos\" (course prerequisites) bullet list in its new\n# order, with three entries dropped, four entries added, and a typo\n# (\"\u00c0lgebra\" -> \"\u00c1lgebra\") fixed along the way.\n\n$doc = $word.ActiveDocument\n$paras = $doc.Paragraphs\n\n$listIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    if ($paras.Item($i).Range.Text.Trim() -eq \"Requisitos\") {\n        $listIndex = $i + 1\n        break\n    }\n}\nif ($listIndex -eq -1) {\n    throw \"Could not find the Requisitos list paragraph\"\n}\n\n$listPara = $paras.Item($listIndex)\n$listRange = $listPara.Range\n# Exclude the trailing paragraph mark from the range so it is never\n# overwritten (this paragraph happens to be the last one in the body).\n$listRange.End = $listRange.End - 1\n$listRange.Text = \"\"\n\n$items = @(\n    \"LOB1268 -  Leitura, Escrita e Comunica\u00e7\u00e3o Cient\u00edfica  (Requisito)\",\n    \"LOB1270 -  Qu\u00edmica Experimental Aplicada  (Requisito)\",\n    \"LOM3081 -  Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos S\u00f3lidos  (Requisito)\",\n    \"LOQ4097 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito)\",\n    \"LOQ4098 -  Fundamentos de Qu\u00edmica para Engenharia II (Requisito)\",\n    \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n    \"LOB1040 -  Laborat\u00f3rio de Eletricidade  (Requisito)\",\n    \"LOB1053 -  F\u00edsica III  (Requisito)\",\n    \"LOB1056 -  Introdu\u00e7\u00e3o aos M\u00e9todos Num\u00e9ricos e Computacionais  (Requisito)\",\n    \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n    \"LOB1011 -  Eletricidade Aplicada  (Requisito)\",\n    \"LOB1012 -  Estat\u00edstica  (Requisito)\",\n    \"LOB1024 -  Mec\u00e2nica  (Requisito)\",\n    \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n    \"LOB1037 -  \u00c1lgebra Linear  (Requisito)\",\n    \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n    \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n    \"LOB1041 -  F\u00edsica Experimental II  (Requisito)\",\n    \"LOB1042 -  F\u00edsica Experimental IV  (Requisito)\",\n    \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n    \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n    \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n    \"LOB1009 -  Leitura e Interpreta\u00e7\u00e3o de Desenho T\u00e9cnico  (Requisito)\",\n    \"LOB1018 -  F\u00edsica I  (Requisito)\",\n    \"LOB1019 -  F\u00edsica II  (Requisito)\",\n    \"LOB1021 -  F\u00edsica IV  (Requisito)\",\n    \"LOQ4233 -  Gest\u00e3o de Neg\u00f3cios  (Requisito)\"\n)\n\n# Re-insert every item in the new order, each followed by a manual line\n# break, matching the original \"one line per requirement\" formatting.\n$insertPos = $listRange.Start\nforeach ($item in $items) {\n    $r = $doc.Range($insertPos, $insertPos)\n    $r.InsertAfter($item)\n    $insertPos = $insertPos + $item.Length\n\n    $br = $doc.Range($insertPos, $insertPos)\n    $br.InsertBreak(6)   # wdLineBreak\n    $insertPos = $insertPos + 1\n}\n"}
